$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the RF (column I) values for rows 19 through 40 from 13.77145454545455
# to 25.77411111111111, reflecting the 2025 data / RF update described in the
# commit message.
$newValue = 25.77411111111111
for ($row = 19; $row -le 40; $row++) {
    $ws.Cells.Item($row, 9).Value = $newValue
}
